$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 30.54296875
$ws.Columns.Item(2).ColumnWidth = 17.54296875

# Cell values
$ws.Range("A1").Value = "tatause02@test.com"
$ws.Range("B1").Value = "Pass@123"
$ws.Range("A2").Value = "seleautouser01@test.com"
$ws.Range("B2").Value = "Pass@123"
$ws.Range("A3").Value = "seleautouser03@test.com"
$ws.Range("B3").Value = "Pass@123"

# Hyperlinks
$ws.Hyperlinks.Add($ws.Range("A1"), "https://demowebshop.tricentis.com/customer/info", "", "", "tatause02@test.com")
$ws.Hyperlinks.Add($ws.Range("B1"), "https://demowebshop.tricentis.com/customer/info", "", "", "Pass@123")
$ws.Hyperlinks.Add($ws.Range("A2"), "https://demowebshop.tricentis.com/customer/info", "", "", "seleautouser01@test.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://demowebshop.tricentis.com/customer/info", "", "", "Pass@123")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://demowebshop.tricentis.com/customer/info", "", "", "seleautouser03@test.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://demowebshop.tricentis.com/customer/info", "", "", "Pass@123")

$ws.Range("B3").Select()
